# Handle format runs in shared string record.
# Renames the sheets, re-selects cells, and adds a "Simple Format" demo
# sheet that shows whole-cell bold/italic plus partial (run-level)
# bold/italic formatting inside a single cell.

$wb = $excel.ActiveWorkbook

# --- Rename the sheets -----------------------------------------------
$wb.Worksheets.Item("Sheet1").Name = "Unformatted"
$wb.Worksheets.Item("Sheet2").Name = "Simple Format"

$ws1 = $wb.Worksheets.Item("Unformatted")
$ws2 = $wb.Worksheets.Item("Simple Format")

# --- Sheet1 ("Unformatted"): move the selection ------------------------
$ws1.Range("G11").Select()

# --- Sheet2 ("Simple Format"): new content -----------------------------

# Row 1: whole-cell bold / whole-cell italic / mixed run formatting
$ws2.Range("A1").Value = "Bold"
$ws2.Range("A1").Font.Bold = $true

$ws2.Range("B1").Value = "Italic"
$ws2.Range("B1").Font.Italic = $true

$ws2.Range("C1").Value = "Part bold and part italic"
$ws2.Range("C1").Characters(6, 4).Font.Bold = $true
$ws2.Range("C1").Characters(10, 10).Font.Bold = $false
$ws2.Range("C1").Characters(20, 6).Font.Italic = $true

# Row 2: the same text shown unformatted, whole-cell bold, whole-cell
# italic, and with mixed run formatting - so the four columns can be
# compared side by side.
$ws2.Range("A2").Value = "Same Text"

$ws2.Range("B2").Value = "Same Text"
$ws2.Range("B2").Font.Bold = $true

$ws2.Range("C2").Value = "Same Text"
$ws2.Range("C2").Font.Italic = $true

$ws2.Range("D2").Value = "Same Text"
$ws2.Range("D2").Characters(1, 4).Font.Bold = $true
$ws2.Range("D2").Characters(5, 1).Font.Bold = $false
$ws2.Range("D2").Characters(6, 4).Font.Italic = $true

# Column widths for the new sheet (A is best-fit on the longest value in
# the column, C is a fixed custom width for the "mixed formatting" demo
# text)
$ws2.Columns.Item(1).ColumnWidth = 9
$ws2.Columns.Item(3).ColumnWidth = 25

# Basic page setup for the new sheet
$ws2.PageSetup.Orientation = 1

# Selection + activation: Sheet2 becomes the active tab with F3 selected.
$ws2.Range("F3").Select()
$ws2.Activate()

Write-Output "done"
